$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add category and title values for row 2
$ws.Range("D2").Value = "cards"
$ws.Range("E2").Value = "Poles Card"

# Add weights: insure, study, job, live
$ws.Range("I2").Value = 5
$ws.Range("J2").Value = 5
$ws.Range("K2").Value = 5
$ws.Range("L2").Value = 10

# Match the author's final selection/view state (scrolled to show column D, cell L2 selected)
$win = $excel.ActiveWindow
try {
    $win.ScrollColumn = 4
    $win.ScrollRow = 2
} catch {
    # best effort only; not all hosts expose window scroll position
}
$ws.Range("L2").Select()
